$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Widen the "zh-cn" and "de-de" columns on the Overview sheet (E, F)
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---- zh-cn sheet ----
# Status column (C) text changed, widen column to fit new text
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527

# Fill in Latest Target File (I), Latest Handback File (J), Latest Handback DateTime (K)
$wsZhCn.Range("I2").Value = "0a31343b-1c44-43a0-ace4-ad00deeb601c.md"
$wsZhCn.Range("I2").Style = "Hyperlink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8325348c6fe67c5946bb640590fcf34982d94c12/e2e/0a31343b-1c44-43a0-ace4-ad00deeb601c.md", [Type]::Missing, [Type]::Missing, "0a31343b-1c44-43a0-ace4-ad00deeb601c.md") | Out-Null
$wsZhCn.Range("J2").Value = "0a31343b-1c44-43a0-ace4-ad00deeb601c.50bf17ec4fc5a347971fb4e82bb0ce2b8849c516.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-18 12:29:14"

$wsZhCn.Range("I3").Value = "ba79d78b-78dd-4a6d-82f1-d7c05090d214.md"
$wsZhCn.Range("I3").Style = "Hyperlink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8325348c6fe67c5946bb640590fcf34982d94c12/e2e/ba79d78b-78dd-4a6d-82f1-d7c05090d214.md", [Type]::Missing, [Type]::Missing, "ba79d78b-78dd-4a6d-82f1-d7c05090d214.md") | Out-Null
$wsZhCn.Range("J3").Value = "ba79d78b-78dd-4a6d-82f1-d7c05090d214.679c12effa595fb4de604b3a64ea4b8f5a2bc4e0.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-18 12:29:14"

$wsZhCn.Columns.Item(9).ColumnWidth = 40
$wsZhCn.Columns.Item(10).ColumnWidth = 40

# ---- de-de sheet ----
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527

$wsDeDe.Range("I2").Value = "0a31343b-1c44-43a0-ace4-ad00deeb601c.md"
$wsDeDe.Range("I2").Style = "Hyperlink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8325348c6fe67c5946bb640590fcf34982d94c12/e2e/0a31343b-1c44-43a0-ace4-ad00deeb601c.md", [Type]::Missing, [Type]::Missing, "0a31343b-1c44-43a0-ace4-ad00deeb601c.md") | Out-Null
$wsDeDe.Range("J2").Value = "0a31343b-1c44-43a0-ace4-ad00deeb601c.50bf17ec4fc5a347971fb4e82bb0ce2b8849c516.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-18 12:29:21"

$wsDeDe.Range("I3").Value = "ba79d78b-78dd-4a6d-82f1-d7c05090d214.md"
$wsDeDe.Range("I3").Style = "Hyperlink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8325348c6fe67c5946bb640590fcf34982d94c12/e2e/ba79d78b-78dd-4a6d-82f1-d7c05090d214.md", [Type]::Missing, [Type]::Missing, "ba79d78b-78dd-4a6d-82f1-d7c05090d214.md") | Out-Null
$wsDeDe.Range("J3").Value = "ba79d78b-78dd-4a6d-82f1-d7c05090d214.679c12effa595fb4de604b3a64ea4b8f5a2bc4e0.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-18 12:29:21"

$wsDeDe.Columns.Item(9).ColumnWidth = 40
$wsDeDe.Columns.Item(10).ColumnWidth = 40

$wb.Save()
